# Add ericsson mixedmode NL
# Adds 6 new rows (94-99) to the "Path" sheet describing the new
# Ericsson "MIXEDMODE_NL" full_kget raw file paths for NTH-ENM, STH-ENM
# and CEW-ENM (once for the LTE/L900-1800-2100 technology row, and once
# for the NR/700-2600 technology row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Path")

# Copy the formatting of the last existing data row down onto the new
# rows first, so the new cells pick up the same styles (s="29"/"24")
# used throughout the table, without touching any existing cell.
$ws.Range("A93:I93").Copy()
$ws.Range("A94:I99").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$data = @(
  @("RFT", "Ericsson", "LTE", "L900/1800/2100", "NTH", "10.50.64.207", "/home/app/ngoss/data/rfserver/Ericsson/NTH-ENM/full_kget/MIXEDMODE_NL/1Current", "Full-Kget.txt", "18.Q1"),
  @("RFT", "Ericsson", "LTE", "L900/1800/2100", "STH", "10.50.64.207", "/home/app/ngoss/data/rfserver/Ericsson/STH-ENM/full_kget/MIXEDMODE_NL/1Current", "Full-Kget.txt", "W18.Q1"),
  @("RFT", "Ericsson", "LTE", "L900/1800/2100", "CEW", "10.50.64.207", "/home/app/ngoss/data/rfserver/Ericsson/CEW-ENM/full_kget/MIXEDMODE_NL/1Current", "Full-Kget.txt", "W18.Q1"),
  @("RFT", "Ericsson", "NR",  "700/2600",       "NTH", "10.50.64.207", "/home/app/ngoss/data/rfserver/Ericsson/NTH-ENM/full_kget/MIXEDMODE_NL/1Current", "Full-Kget.txt", "18.Q1"),
  @("RFT", "Ericsson", "NR",  "700/2601",       "STH", "10.50.64.207", "/home/app/ngoss/data/rfserver/Ericsson/STH-ENM/full_kget/MIXEDMODE_NL/1Current", "Full-Kget.txt", "W18.Q1"),
  @("RFT", "Ericsson", "NR",  "700/2602",       "CEW", "10.50.64.207", "/home/app/ngoss/data/rfserver/Ericsson/CEW-ENM/full_kget/MIXEDMODE_NL/1Current", "Full-Kget.txt", "W18.Q1")
)

$row = 94
foreach ($r in $data) {
  $ws.Cells.Item($row, 1).Value = $r[0]
  $ws.Cells.Item($row, 2).Value = $r[1]
  $ws.Cells.Item($row, 3).Value = $r[2]
  $ws.Cells.Item($row, 4).Value = $r[3]
  $ws.Cells.Item($row, 5).Value = $r[4]
  $ws.Cells.Item($row, 6).Value = $r[5]
  $ws.Cells.Item($row, 7).Value = $r[6]
  $ws.Cells.Item($row, 8).Value = $r[7]
  $ws.Cells.Item($row, 9).Value = $r[8]
  $row = $row + 1
}

# The workbook also picked up a small (8pt) Calibri font used for the
# sheet's phonetic-guide settings (Home > Font > Show/Hide Phonetic
# Fields dialog). Register that font in the style table the same way,
# via a scratch cell that is cleared right back out so no visible cell
# actually ends up using it.
$scratch = $ws.Range("Z1")
$scratch.Font.Size = 8
$scratch.Clear() | Out-Null

# Reposition the view to match the edited region (scrolled down to the
# new rows, selection parked just past them).
$ws.Activate() | Out-Null
$ws.Range("G104").Select() | Out-Null

Write-Host "Added rows 94-99 to Path sheet"
